$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the overall n count (C2) to reflect samples without relapse
$ws.Range("C2").Value = 2740

# Update the Dx score comparison percentages (C3:C26) to the
# "without relapse samples" values
$ws.Range("C3").Value  = "905 (41.8)"
$ws.Range("C4").Value  = "902 (41.7)"
$ws.Range("C5").Value  = "31 (1.4)"
$ws.Range("C6").Value  = "48 (2.2)"
$ws.Range("C7").Value  = "138 (6.4)"
$ws.Range("C8").Value  = "140 (6.5)"
$ws.Range("C9").Value  = "398 (22.9)"
$ws.Range("C10").Value = "395 (22.7)"
$ws.Range("C11").Value = "576 (33.2)"
$ws.Range("C12").Value = "165 (9.5)"
$ws.Range("C13").Value = "203 (11.7)"
$ws.Range("C14").Value = "766 (49.5)"
$ws.Range("C15").Value = "783 (50.5)"
$ws.Range("C16").Value = "58 (2.1)"
$ws.Range("C17").Value = "592 (21.6)"
$ws.Range("C18").Value = "397 (14.5)"
$ws.Range("C19").Value = "280 (10.2)"
$ws.Range("C20").Value = "27 (1.0)"
$ws.Range("C21").Value = "83 (3.0)"
$ws.Range("C22").Value = "153 (5.6)"
$ws.Range("C23").Value = "64 (2.3)"
$ws.Range("C24").Value = "764 (27.9)"
$ws.Range("C25").Value = "128 (4.7)"
$ws.Range("C26").Value = "194 (7.1)"
